$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A99").Value = "GRT-USD"
